$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5286
$ws.Range("J3").Value = 5600
$ws.Range("J4").Value = 1243
$ws.Range("J5").Value = 439
$ws.Range("J6").Value = 7006
$ws.Range("J7").Value = 19574

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 48

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J2").Value = 51
$ws.Range("J3").Value = 36
$ws.Range("J6").Value = 165
$ws.Range("J7").Value = 260

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J6").Value = 411
$ws.Range("J7").Value = 1237

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J4").Value = 10
$ws.Range("J6").Value = 111

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J2").Value = 222
$ws.Range("J3").Value = 292
$ws.Range("J4").Value = 38
$ws.Range("J6").Value = 304
$ws.Range("J7").Value = 896

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("J3").Value = 98
$ws.Range("J5").Value = 9
$ws.Range("J7").Value = 291

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 211
$ws.Range("J6").Value = 175
$ws.Range("J7").Value = 612

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 145
$ws.Range("J6").Value = 180
$ws.Range("J7").Value = 505

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("J4").Value = 14
$ws.Range("J7").Value = 306

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J4").Value = 74
$ws.Range("J7").Value = 574
$ws.Range("J8").Value = 1237
$ws.Range("J13").Value = 24
$ws.Range("J15").Value = 213
$ws.Range("J18").Value = 165
$ws.Range("J19").Value = 561
$ws.Range("J20").Value = 409
$ws.Range("J25").Value = 98
$ws.Range("J27").Value = 116
$ws.Range("J29").Value = 1094
$ws.Range("J33").Value = 896
$ws.Range("J34").Value = 91
$ws.Range("J36").Value = 271
$ws.Range("J37").Value = 612
$ws.Range("J40").Value = 45
$ws.Range("J42").Value = 796
$ws.Range("J47").Value = 147
$ws.Range("J49").Value = 131
$ws.Range("J50").Value = 122
$ws.Range("J52").Value = 493
$ws.Range("J53").Value = 260
$ws.Range("J54").Value = 377
$ws.Range("J55").Value = 256
$ws.Range("J59").Value = 25
$ws.Range("J60").Value = 122
$ws.Range("J63").Value = 72
$ws.Range("J64").Value = 132
$ws.Range("J65").Value = 505
$ws.Range("J66").Value = 62
$ws.Range("J67").Value = 752
$ws.Range("J69").Value = 48
$ws.Range("J73").Value = 179
$ws.Range("J76").Value = 281
$ws.Range("J79").Value = 560
$ws.Range("J82").Value = 26
$ws.Range("J84").Value = 167
$ws.Range("J86").Value = 120
$ws.Range("J88").Value = 213
$ws.Range("J89").Value = 251
$ws.Range("J91").Value = 218
$ws.Range("J93").Value = 86
$ws.Range("J94").Value = 193
$ws.Range("J95").Value = 291
$ws.Range("J98").Value = 131
$ws.Range("J99").Value = 306
$ws.Range("J101").Value = 19574

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 188
$ws.Range("J3").Value = 290
$ws.Range("J7").Value = 752

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J2").Value = 54
$ws.Range("J3").Value = 52
$ws.Range("J7").Value = 167

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J6").Value = 76
$ws.Range("J7").Value = 131

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J6").Value = 180
$ws.Range("J7").Value = 377

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 374
$ws.Range("J6").Value = 288
$ws.Range("J7").Value = 1094

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J3").Value = 164
$ws.Range("J6").Value = 206
$ws.Range("J7").Value = 561

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J2").Value = 43
$ws.Range("J3").Value = 57
$ws.Range("J6").Value = 155
$ws.Range("J7").Value = 281

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J4").Value = 37
$ws.Range("J7").Value = 796

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("J3").Value = 8
$ws.Range("J6").Value = 24

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 123
$ws.Range("J7").Value = 256

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J2").Value = 66
$ws.Range("J7").Value = 218

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 197
$ws.Range("J6").Value = 157
$ws.Range("J7").Value = 560

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J2").Value = 37
$ws.Range("J7").Value = 132

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 112
$ws.Range("J3").Value = 140
$ws.Range("J7").Value = 409

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J2").Value = 43
$ws.Range("J7").Value = 165

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J3").Value = 87
$ws.Range("J6").Value = 78
$ws.Range("J7").Value = 271

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("J3").Value = 26
$ws.Range("J7").Value = 86

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 179
$ws.Range("J3").Value = 174
$ws.Range("J5").Value = 16
$ws.Range("J6").Value = 182
$ws.Range("J7").Value = 574

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J5").Value = 3
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J3").Value = 39
$ws.Range("J6").Value = 104
$ws.Range("J7").Value = 193

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J6").Value = 16
$ws.Range("J7").Value = 98

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J2").Value = 34
$ws.Range("J7").Value = 147

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J6").Value = 90
$ws.Range("J7").Value = 213

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J2").Value = 23
$ws.Range("J6").Value = 81
$ws.Range("J7").Value = 131

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 122

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J2").Value = 63
$ws.Range("J7").Value = 179

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("J6").Value = 6
$ws.Range("J7").Value = 25

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J6").Value = 95
$ws.Range("J7").Value = 213

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("J2").Value = 81
$ws.Range("J3").Value = 69
$ws.Range("J7").Value = 251

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J2").Value = 34
$ws.Range("J7").Value = 116

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J2").Value = 16
$ws.Range("J7").Value = 120

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J2").Value = 44
$ws.Range("J7").Value = 122

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("J5").Value = 17
$ws.Range("J6").Value = 26

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("J2").Value = 18
$ws.Range("J7").Value = 45

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 119
$ws.Range("J7").Value = 493

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J2").Value = 26
$ws.Range("J7").Value = 74
